$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 50,4

$data[0,0] = 0.07109836488962173
$data[0,1] = 0.9812811017036438
$data[0,2] = 0.01564565859735012
$data[0,3] = 0.9975326657295227
$data[1,0] = 0.009075576439499855
$data[1,1] = 0.9984514117240906
$data[1,2] = 0.01309419330209494
$data[1,3] = 0.997738242149353
$data[2,0] = 0.004349478986114264
$data[2,1] = 0.9989740252494812
$data[2,2] = 0.004659167490899563
$data[2,3] = 0.998766303062439
$data[3,0] = 0.002166710095480084
$data[3,1] = 0.9992644190788269
$data[3,2] = 0.002583120251074433
$data[3,3] = 0.9995887875556946
$data[4,0] = 0.001694503938779235
$data[4,1] = 0.9996128678321838
$data[4,2] = 0.006677067838609219
$data[4,3] = 0.9990747570991516
$data[5,0] = 0.001663376926444471
$data[5,1] = 0.9995934963226318
$data[5,2] = 0.001893842127174139
$data[5,3] = 0.9997943639755249
$data[6,0] = 0.001703369431197643
$data[6,1] = 0.9996128678321838
$data[6,2] = 0.002304930472746491
$data[6,3] = 0.9997943639755249
$data[7,0] = 0.000576861493755132
$data[7,1] = 0.9998645186424255
$data[7,2] = 0.001988050295040011
$data[7,3] = 0.9997943639755249
$data[8,0] = 0.0006628966075368226
$data[8,1] = 0.9998451471328735
$data[8,2] = 0.002106309169903398
$data[8,3] = 0.9998971819877625
$data[9,0] = 0.001153858262114227
$data[9,1] = 0.9997096061706543
$data[9,2] = 0.001744475099258125
$data[9,3] = 0.9998971819877625
$data[10,0] = 0.0003597549512051046
$data[10,1] = 0.9999225735664368
$data[10,2] = 0.002221499336883426
$data[10,3] = 0.9998971819877625
$data[11,0] = 0.0006011043442413211
$data[11,1] = 0.9998645186424255
$data[11,2] = 0.002131890971213579
$data[11,3] = 0.9998971819877625
$data[12,0] = 0.0009197793551720679
$data[12,1] = 0.9998838305473328
$data[12,2] = 0.002157262759283185
$data[12,3] = 0.9998971819877625
$data[13,0] = 0.0002216671855421737
$data[13,1] = 0.999961256980896
$data[13,2] = 0.002732397057116032
$data[13,3] = 0.9997943639755249
$data[14,0] = 0.001313347369432449
$data[14,1] = 0.9998064041137695
$data[14,2] = 0.00197906862013042
$data[14,3] = 0.9998971819877625
$data[15,0] = 0.0005235631833784282
$data[15,1] = 0.9998451471328735
$data[15,2] = 0.002374853938817978
$data[15,3] = 0.9998971819877625
$data[16,0] = 0.0004051274445373565
$data[16,1] = 0.9999032020568848
$data[16,2] = 0.00231908424757421
$data[16,3] = 0.9998971819877625
$data[17,0] = 0.000418732495745644
$data[17,1] = 0.9998451471328735
$data[17,2] = 0.003000837750732899
$data[17,3] = 0.9997943639755249
$data[18,0] = 0.000467491481686011
$data[18,1] = 0.9999032020568848
$data[18,2] = 0.002852467121556401
$data[18,3] = 0.9997943639755249
$data[19,0] = 0.0008965424494817853
$data[19,1] = 0.9998257756233215
$data[19,2] = 0.002433480462059379
$data[19,3] = 0.9998971819877625
$data[20,0] = 0.0003604766679927707
$data[20,1] = 0.9999032020568848
$data[20,2] = 0.002840063767507672
$data[20,3] = 0.9998971819877625
$data[21,0] = 0.0004379321180749685
$data[21,1] = 0.9998257756233215
$data[21,2] = 0.003185037523508072
$data[21,3] = 0.9998971819877625
$data[22,0] = 0.0002963498118333519
$data[22,1] = 0.9999032020568848
$data[22,2] = 0.003401753725484014
$data[22,3] = 0.9998971819877625
$data[23,0] = 0.00005295739174471237
$data[23,1] = 1
$data[23,2] = 0.003776739817112684
$data[23,3] = 0.9998971819877625
$data[24,0] = 0.00004238043402438052
$data[24,1] = 1
$data[24,2] = 0.003912874031811953
$data[24,3] = 0.9998971819877625
$data[25,0] = 0.0007558464421890676
$data[25,1] = 0.9999032020568848
$data[25,2] = 0.003969909157603979
$data[25,3] = 0.9997943639755249
$data[26,0] = 0.0003738219093065709
$data[26,1] = 0.9998645186424255
$data[26,2] = 0.004174177069216967
$data[26,3] = 0.9998971819877625
$data[27,0] = 0.0005741852801293135
$data[27,1] = 0.9998451471328735
$data[27,2] = 0.003541883081197739
$data[27,3] = 0.9998971819877625
$data[28,0] = 0.0004834337451029569
$data[28,1] = 0.9999032020568848
$data[28,2] = 0.003674410283565521
$data[28,3] = 0.9998971819877625
$data[29,0] = 0.0003383196890354156
$data[29,1] = 0.9999032020568848
$data[29,2] = 0.003607323160395026
$data[29,3] = 0.9998971819877625
$data[30,0] = 0.0003489677328616381
$data[30,1] = 0.9999225735664368
$data[30,2] = 0.003572541056200862
$data[30,3] = 0.9998971819877625
$data[31,0] = 0.0005261022015474737
$data[31,1] = 0.9999032020568848
$data[31,2] = 0.003298641415312886
$data[31,3] = 0.9998971819877625
$data[32,0] = 0.0001544979750178754
$data[32,1] = 0.9999225735664368
$data[32,2] = 0.003825076157227159
$data[32,3] = 0.9998971819877625
$data[33,0] = 0.0000115028788059135
$data[33,1] = 1
$data[33,2] = 0.003960392437875271
$data[33,3] = 0.9998971819877625
$data[34,0] = 0.0003065828641410917
$data[34,1] = 0.9999225735664368
$data[34,2] = 0.004265422467142344
$data[34,3] = 0.9998971819877625
$data[35,0] = 0.0004491801664698869
$data[35,1] = 0.9998645186424255
$data[35,2] = 0.003699644468724728
$data[35,3] = 0.9998971819877625
$data[36,0] = 0.0002785635297186673
$data[36,1] = 0.9999032020568848
$data[36,2] = 0.003723705653101206
$data[36,3] = 0.9998971819877625
$data[37,0] = 0.0004806426004506648
$data[37,1] = 0.9999419450759888
$data[37,2] = 0.00336453365162015
$data[37,3] = 0.9998971819877625
$data[38,0] = 0.00002330929419258609
$data[38,1] = 1
$data[38,2] = 0.004090424627065659
$data[38,3] = 0.9998971819877625
$data[39,0] = 0.0005164096364751458
$data[39,1] = 0.9999225735664368
$data[39,2] = 0.00367808504961431
$data[39,3] = 0.9998971819877625
$data[40,0] = 0.000139405092340894
$data[40,1] = 0.9999419450759888
$data[40,2] = 0.003960953559726477
$data[40,3] = 0.9998971819877625
$data[41,0] = 0.0002511959464754909
$data[41,1] = 0.999980628490448
$data[41,2] = 0.003251246176660061
$data[41,3] = 0.9998971819877625
$data[42,0] = 0.00005777502155979164
$data[42,1] = 0.999980628490448
$data[42,2] = 0.003847262356430292
$data[42,3] = 0.9998971819877625
$data[43,0] = 0.0003823455772362649
$data[43,1] = 0.9999419450759888
$data[43,2] = 0.004236360546201468
$data[43,3] = 0.9998971819877625
$data[44,0] = 0.000301788590149954
$data[44,1] = 0.9999032020568848
$data[44,2] = 0.004116904456168413
$data[44,3] = 0.9998971819877625
$data[45,0] = 0.00008739442273508757
$data[45,1] = 0.999961256980896
$data[45,2] = 0.004202902317047119
$data[45,3] = 0.9998971819877625
$data[46,0] = 0.000006159298663988011
$data[46,1] = 1
$data[46,2] = 0.004412871785461903
$data[46,3] = 0.9998971819877625
$data[47,0] = 0.0002314865414518863
$data[47,1] = 0.999961256980896
$data[47,2] = 0.004467417486011982
$data[47,3] = 0.9998971819877625
$data[48,0] = 0.0009304629638791084
$data[48,1] = 0.9998645186424255
$data[48,2] = 0.004273765720427036
$data[48,3] = 0.9998971819877625
$data[49,0] = 0.00001377763601340121
$data[49,1] = 1
$data[49,2] = 0.004614434670656919
$data[49,3] = 0.9998971819877625

$ws.Range("A2:D51").Value = $data
